$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "time_taken" header in F1, copying the header style from E1
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the time_taken values for each data row
$ws.Range("F2").Value = "2021-10-05 13:41:04.873901"
$ws.Range("F3").Value = "2021-10-05 13:41:04.873912"
$ws.Range("F4").Value = "2021-10-05 13:41:04.873915"
$ws.Range("F5").Value = "2021-10-05 13:41:04.873918"
$ws.Range("F6").Value = "2021-10-05 13:41:04.873921"
$ws.Range("F7").Value = "2021-10-05 13:41:04.873924"
$ws.Range("F8").Value = "2021-10-05 13:41:04.873926"
$ws.Range("F9").Value = "2021-10-05 13:41:04.873929"
$ws.Range("F10").Value = "2021-10-05 13:41:04.873932"
$ws.Range("F11").Value = "2021-10-05 13:41:04.873934"
$ws.Range("F12").Value = "2021-10-05 13:41:04.873937"
$ws.Range("F13").Value = "2021-10-05 13:41:04.873939"
$ws.Range("F14").Value = "2021-10-05 13:41:04.873942"
$ws.Range("F15").Value = "2021-10-05 13:41:04.873944"
$ws.Range("F16").Value = "2021-10-05 13:41:04.873947"
$ws.Range("F17").Value = "2021-10-05 13:41:04.873949"
$ws.Range("F18").Value = "2021-10-05 13:41:04.873952"
$ws.Range("F19").Value = "2021-10-05 13:41:04.873955"

Write-Host "Added time_taken column F1:F19"
